# Apply cryptocurrency price/volume updates (and a few row re-orderings)
# as captured by the commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.741.65"
$ws.Range("E2").Value = "  +1.49%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.881.85"
$ws.Range("E3").Value = "  +1.40%  "
# Row 4
$ws.Range("E4").Value = "  +0.44%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.12"
$ws.Range("E5").Value = "  +2.03%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.42%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4706"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3934"
$ws.Range("E8").Value = "  +0.82%  "
# Row 9
$ws.Range("E9").Value = "  -0.84%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08059"
$ws.Range("E10").Value = "  +1.64%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.027"
$ws.Range("E11").Value = "  +1.59%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.19"
$ws.Range("E12").Value = "  +3.54%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.886.24"
$ws.Range("E13").Value = "  +1.78%  "
# Row 14
$ws.Range("E14").Value = "  +0.94%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.142"
$ws.Range("E15").Value = "  -0.36%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("E16").Value = "  +0.50%  "
# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06716"
$ws.Range("E17").Value = "  +1.57%  "
# Row 18
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.22"
$ws.Range("E18").Value = "  +1.36%  "
# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001048"
$ws.Range("E19").Value = "  +1.67%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.34"
$ws.Range("E20").Value = "  +0.59%  "
# Row 21
$ws.Range("E21").Value = "  +0.30%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.544"
$ws.Range("E22").Value = "  +0.71%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.745.20"
$ws.Range("E23").Value = "  +1.50%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.03"
$ws.Range("E24").Value = "  +1.46%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.313"
$ws.Range("E25").Value = "  +0.93%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.109.59"
$ws.Range("E26").Value = "  +1.66%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.97"
$ws.Range("E27").Value = "  +3.90%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.19"
$ws.Range("E28").Value = "  +1.20%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.109"
$ws.Range("E29").Value = "  +2.08%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.588"
$ws.Range("E30").Value = "  +2.33%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.91"
$ws.Range("E31").Value = "  +0.67%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9875"
$ws.Range("E32").Value = "  +3.92%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09478"
$ws.Range("E33").Value = "  +1.30%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.452"
$ws.Range("E34").Value = "  +0.08%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.620"
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.359"
$ws.Range("E36").Value = "  +1.86%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06138"
$ws.Range("E37").Value = "  +1.76%  "
# Row 38
$ws.Range("E38").Value = "  +1.74%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.232"
$ws.Range("E39").Value = "  +1.00%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.138"
$ws.Range("E40").Value = "  +1.05%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6009"
$ws.Range("E41").Value = "  +1.37%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1899"
$ws.Range("E42").Value = "  +0.64%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.32"
$ws.Range("E43").Value = "  +1.52%  "
# Row 44
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.263"
$ws.Range("E44").Value = "  -1.43%  "
# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5730"
$ws.Range("E45").Value = "  +1.98%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.21"
$ws.Range("E46").Value = "  +0.64%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.947"
$ws.Range("E47").Value = "  +1.37%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.391"
$ws.Range("E48").Value = "  +0.29%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06905"
$ws.Range("E49").Value = "  +2.46%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.34"
$ws.Range("E50").Value = "  +5.60%  "
# Row 51
$ws.Range("E51").Value = "  +1.79%  "
